$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (BOM_BUCKET 3)
$ws.Range("B2").Value = 18749565.12
$ws.Range("C2").Value = 493242
$ws.Range("D2").Value = 2.63

# Update row 3 (BOM_BUCKET 4)
$ws.Range("B3").Value = 15513891.85
$ws.Range("C3").Value = 171573
$ws.Range("D3").Value = 1.11

# Update row 4 (BOM_BUCKET 5)
$ws.Range("B4").Value = 9530315.050000001
$ws.Range("C4").Value = 32272
$ws.Range("D4").Value = 0.34
